$d = $word.ActiveDocument

# 1. Add sentence about flow diagram after the first-exercise intro sentence.
$d.Content.Find.Execute(
    "As a first exercise we are going to run the simple SEIR model, as seen in practical 2, in R.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As a first exercise we are going to run the simple SEIR model, as seen in practical 2, in R. See practical 2. for the flow diagram.",
    2)

# 2. Fix typo "sizr" -> "size" in the interactivity comment.
$d.Content.Find.Execute(
    "## Interactivity allows plot zooming and gives a tool tip providing the population sizr at any point.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## Interactivity allows plot zooming and gives a tool tip providing the population size at any point.",
    2)

# 3. Add missing "2" after "See practical" in the SHLIR intro paragraph.
$d.Content.Find.Execute(
    "Now we are going to implement the SHLIR model from practical 2 and try to reproduce some of the behaviour observed using the interactive interface. See practical for details on this model.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Now we are going to implement the SHLIR model from practical 2 and try to reproduce some of the behaviour observed using the interactive interface. See practical 2 for details on this model.",
    2)

# 4. Only the FIRST "Can you alter the parameters..." bullet (under the SEIR
#    explore section, numId 1003) should mention the SHLIR model explicitly;
#    a second, identical bullet later in the document (under the SHLIR
#    explore section, numId 1004) must stay unchanged. Find the first
#    paragraph matching that exact text and replace just its range.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd().TrimEnd([char]13) -eq "Can you alter the parameters defined above to answer the questions for this model from practical 2?") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $rng = $target.Range
    $rng.Find.Execute(
        "Can you alter the parameters defined above to answer the questions for this model from practical 2?",
        $true, $false, $false, $false, $false, $true, 1, $false,
        "Can you alter the parameters defined above to answer the questions for the SHLIR model from practical 2?",
        2)
}
